# Build site at 2023-01-09 16:18:13 UTC
# Update the LOM3204 discipline sheet:
#  - Ativacao (activation date) changes from 01/01/2016 to 01/01/2023
#    (this value is also (mistakenly) reused by the "Programa resumido:" row,
#    i.e. cells B13/C13 point at the very same text as B8/C8)
#  - New English translations are filled in for:
#      Objectives:      -> long description (row 11)
#      Short syllabus:  -> short description (row 14)
#      Syllabus:        -> long description (row 16)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update the activation date text in B8/C8 and B13/C13.
#    A plain assignment of "01/01/2023" gets auto-recognised as a date by
#    the spreadsheet engine, which would turn the cell into a numeric date
#    value instead of the original text string. To keep it as literal
#    text (matching the original shared-string cell type) we build the
#    text in a scratch cell using a formula (which always yields a text
#    result) and then paste only the *value* of that scratch cell onto
#    the target cells, preserving their existing formatting.
# ---------------------------------------------------------------------
$ws.Range("Z1").Formula = "=T(""01/01/2023"")"
$ws.Range("Z1").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

# ---------------------------------------------------------------------
# 2) Row 11 "Objectives:" -> add the English objectives text in B11/C11.
#    Copy formatting from the row above (B10/C10) so the new cells get
#    the same styles (wrap text, normal / red font) as the rest of the
#    table before setting their text.
# ---------------------------------------------------------------------
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$objectivesText = "Develop knowledge in order to make the student capable of correctly interpreting the technical drawing, knowing the methodologies and tools used in the industry, giving subsidies so that they can execute, interact and modify drawings and projects throughout their professional life."
$ws.Range("B11").Value = $objectivesText
$ws.Range("C11").Value = $objectivesText

# ---------------------------------------------------------------------
# 3) Row 14 "Short syllabus:" -> add the English short syllabus text in
#    B14/C14.
# ---------------------------------------------------------------------
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$shortSyllabusText = "Context of the technical drawing in the industry, main tools and techniques used in drawings for the elaboration of projects. Introduction to computer-aided design (CAD)."
$ws.Range("B14").Value = $shortSyllabusText
$ws.Range("C14").Value = $shortSyllabusText

# ---------------------------------------------------------------------
# 4) Row 16 "Syllabus:" -> add the English syllabus text in B16/C16.
# ---------------------------------------------------------------------
$ws.Range("B10").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$syllabusText = "Technical drawing standards. Technical terminology and materials for drawing. Perspective representation. Orthogonal design. Scaling and scaling. Cut and section. Auxiliary view and details. Geometric tolerances. Representation of machine elements. Use of software for technical design. Computer-aided design in three dimensions (Solid Modeling). Computer-aided design in two dimensions."
$ws.Range("B16").Value = $syllabusText
$ws.Range("C16").Value = $syllabusText

$wb.Save()
